$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $text = [string]$val
    if ($text -eq "") {
        continue
    }

    $parts = @($text -split ", ")

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) {
                $rest += $p
            }
        }
        $newParts = @("System") + $rest
        $newText = $newParts -join ", "
        $cell.Value = $newText
    }
}
